# Update the "取得日時" (retrieved timestamp) column (A) for all data rows
# on the "ランサーズ" sheet from the old scrape time to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-16 01:40:08"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
